$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark from the very first paragraph
#    (was wrapping nothing, right after the "...我描述的有差异" run).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Expand / correct the "streaming" bullet's text.
# ---------------------------------------------------------------------------
$old = "在streaming同时开启的情况下，前者与后者分别是packing不开启与开启的训练参数对比，完全一致"
$new = "在streaming同时开启，验证了使用alpaca数据对象格式，jsonl存储方式的情况下，前者与后者分别是packing不开启与开启的训练参数对比，完全一致。"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ---------------------------------------------------------------------------
# 3) Re-add the "_GoBack" bookmark, now collapsed at the end of that same
#    "streaming" paragraph (after the run, before the paragraph mark).
#
#    A collapsed Range sitting exactly one character before a paragraph's
#    end confuses Bookmarks.Add in this host, so we work around it: drop a
#    throw-away character at that spot, wrap the bookmark around it, then
#    delete the character again -- leaving the bookmark correctly collapsed
#    at the original position.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*packing*" -and $p.Range.Text -like "*streaming*") {
        $target = $p
    }
}

$endPos = $target.Range.End - 1
$d.Range($endPos, $endPos).InsertAfter("Z")
$wrap = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $wrap)
$d.Range($endPos, $endPos + 1).Delete()

# ---------------------------------------------------------------------------
# 4) A handful of "numId=0" (explicitly de-listed) paragraphs are missing an
#    explicit <w:ilvl val="0"/> alongside the <w:numId val="0"/>. Re-apply a
#    list level then strip the numbering back off -- that leaves the ilvl
#    element behind while restoring numId to 0. These are the picture-only
#    paragraphs right after: "...正确的jsonl..." / the streaming bullet /
#    the ParquetConvertInstruct bullet / the max_step explanation / the
#    "的设置，max_step..." bullet.
# ---------------------------------------------------------------------------
foreach ($i in 17, 21, 22, 24, 25) {
    $p = $d.Paragraphs($i)
    $p.Range.ListFormat.ListLevelNumber = 1
    $p.Range.ListFormat.RemoveNumbers()
}
